$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.851.76'
Set-TextValue $ws.Range("E2") '  -3.72%  '
Set-TextValue $ws.Range("D3") '3.498.04'
Set-TextValue $ws.Range("E3") '  -4.30%  '
Set-TextValue $ws.Range("E4") '  -0.02%  '
Set-TextValue $ws.Range("D5") '577.99'
Set-TextValue $ws.Range("E5") '  -1.69%  '
Set-TextValue $ws.Range("D6") '176.65'
Set-TextValue $ws.Range("E6") '  -2.07%  '
Set-TextValue $ws.Range("D7") '0.622'
Set-TextValue $ws.Range("E7") '  -0.22%  '
Set-TextValue $ws.Range("D8") '3.496.84'
Set-TextValue $ws.Range("E8") '  -4.20%  '
Set-TextValue $ws.Range("D9") '1.00'
Set-TextValue $ws.Range("E9") '  -0.01%  '
Set-TextValue $ws.Range("D10") '0.189'
Set-TextValue $ws.Range("E10") '  -6.96%  '
Set-TextValue $ws.Range("E11") '  +7.72%  '
Set-TextValue $ws.Range("E12") '  -1.21%  '
Set-TextValue $ws.Range("D13") '47.31'
Set-TextValue $ws.Range("E13") '  -5.04%  '
Set-TextValue $ws.Range("D14") '0.0000277'
Set-TextValue $ws.Range("E14") '  -3.22%  '
Set-TextValue $ws.Range("D15") '685.44'
Set-TextValue $ws.Range("E15") '  +0.64%  '
Set-TextValue $ws.Range("D16") '8.88'
Set-TextValue $ws.Range("E16") '  -1.48%  '
Set-TextValue $ws.Range("D17") '4.054.73'
Set-TextValue $ws.Range("E17") '  -4.32%  '
Set-TextValue $ws.Range("D18") '68.828.03'
Set-TextValue $ws.Range("E18") '  -3.86%  '
Set-TextValue $ws.Range("D19") '3.495.75'
Set-TextValue $ws.Range("E19") '  -5.46%  '
Set-TextValue $ws.Range("E20") '  -1.57%  '
Set-TextValue $ws.Range("D21") '17.53'
Set-TextValue $ws.Range("E21") '  -3.79%  '
Set-TextValue $ws.Range("D22") '11.18'
Set-TextValue $ws.Range("E22") '  -3.92%  '
Set-TextValue $ws.Range("E23") '  -4.04%  '
Set-TextValue $ws.Range("D24") '16.37'
Set-TextValue $ws.Range("E24") '  -8.28%  '
Set-TextValue $ws.Range("D25") '98.15'
Set-TextValue $ws.Range("E25") '  -4.98%  '
Set-TextValue $ws.Range("E26") '  -4.50%  '
Set-TextValue $ws.Range("E27") '  +0.03%  '
Set-TextValue $ws.Range("E28") '  -6.64%  '
Set-TextValue $ws.Range("D29") '9.42'
Set-TextValue $ws.Range("E29") '  -7.78%  '
Set-TextValue $ws.Range("D30") '33.01'
Set-TextValue $ws.Range("D31") '8.76'
Set-TextValue $ws.Range("E31") '  -4.64%  '
Set-TextValue $ws.Range("E32") '  -7.46%  '
Set-TextValue $ws.Range("D33") '7.35'
Set-TextValue $ws.Range("E33") '  -1.36%  '
Set-TextValue $ws.Range("D34") '1.36'
Set-TextValue $ws.Range("E34") '  -5.98%  '
Set-TextValue $ws.Range("D35") '568.88'
Set-TextValue $ws.Range("E35") '  -2.03%  '
Set-TextValue $ws.Range("D36") '3.65'
Set-TextValue $ws.Range("E36") '  -13.45%  '
Set-TextValue $ws.Range("D37") '10.96'
Set-TextValue $ws.Range("E37") '  -3.35%  '
Set-TextValue $ws.Range("E38") '  -3.13%  '
Set-TextValue $ws.Range("D39") '56.86'
Set-TextValue $ws.Range("E39") '  -4.51%  '
Set-TextValue $ws.Range("D40") '0.999'
Set-TextValue $ws.Range("E40") '  +0.09%  '
Set-TextValue $ws.Range("D41") '0.0441'
Set-TextValue $ws.Range("E41") '  -4.69%  '
Set-TextValue $ws.Range("E42") '  -4.22%  '
Set-TextValue $ws.Range("E43") '  -2.95%  '
Set-TextValue $ws.Range("D44") '3.424.84'
Set-TextValue $ws.Range("E44") '  -8.80%  '
Set-TextValue $ws.Range("D45") '33.55'
Set-TextValue $ws.Range("E45") '  -5.81%  '
Set-TextValue $ws.Range("E46") '  -8.03%  '
Set-TextValue $ws.Range("D47") '2.93'
Set-TextValue $ws.Range("E47") '  +4.64%  '
Set-TextValue $ws.Range("E48") '  -7.08%  '
Set-TextValue $ws.Range("E49") '  -0.35%  '
Set-TextValue $ws.Range("D50") '134.35'
Set-TextValue $ws.Range("E50") '  +0.21%  '
Set-TextValue $ws.Range("E51") '  -0.47%  '
